# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on several Leve rows
# across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 8500
$ws.Range("J10").Value = 8500
$ws.Range("L10").Value = 8500
$ws.Range("N10").Value = -9086

$ws.Range("H76").Value = 3311.0527
$ws.Range("I76").Value = 2995.3333
$ws.Range("J76").Value = 4495
$ws.Range("K76").Value = 2995.3333
$ws.Range("L76").Value = 4495
$ws.Range("M76").Value = -2680.3333
$ws.Range("N76").Value = -5125

$ws.Range("H79").Value = 3311.0527
$ws.Range("I79").Value = 2995.3333
$ws.Range("J79").Value = 4495
$ws.Range("K79").Value = 2995.3333
$ws.Range("L79").Value = 4495
$ws.Range("M79").Value = -1903.3333
$ws.Range("N79").Value = -6679

$ws.Range("H116").Value = 2506.1765
$ws.Range("I116").Value = 1949.7
$ws.Range("J116").Value = 3301.1428
$ws.Range("K116").Value = 1949.7
$ws.Range("L116").Value = 3301.1428
$ws.Range("M116").Value = 1492.3
$ws.Range("N116").Value = -10185.1428

$ws.Range("H138").Value = 4432.8726
$ws.Range("J138").Value = 5021.1562
$ws.Range("L138").Value = 15063.4686
$ws.Range("N138").Value = -25343.4686

$ws.Range("H141").Value = 3335.9524
$ws.Range("I141").Value = 2210
$ws.Range("J141").Value = 6150.8335
$ws.Range("K141").Value = 6630
$ws.Range("L141").Value = 18452.5005
$ws.Range("M141").Value = -1450
$ws.Range("N141").Value = -28812.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 99
$ws.Range("I3").Value = 99
$ws.Range("K3").Value = 99
$ws.Range("M3").Value = 16

$ws.Range("H61").Value = 1649.9387
$ws.Range("I61").Value = 1436.3448
$ws.Range("J61").Value = 1959.65
$ws.Range("K61").Value = 1436.3448
$ws.Range("L61").Value = 1959.65
$ws.Range("M61").Value = -1224.3448
$ws.Range("N61").Value = -2383.65

$ws.Range("H132").Value = 5779.717
$ws.Range("I132").Value = 5286.793
$ws.Range("K132").Value = 15860.379
$ws.Range("M132").Value = -13330.379

$ws.Range("H136").Value = 1649.9387
$ws.Range("I136").Value = 1436.3448
$ws.Range("J136").Value = 1959.65
$ws.Range("K136").Value = 4309.0344
$ws.Range("L136").Value = 5878.950000000001
$ws.Range("M136").Value = -1759.0344
$ws.Range("N136").Value = -10978.95

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 395.41666
$ws.Range("I22").Value = 287.22223
$ws.Range("J22").Value = 720
$ws.Range("K22").Value = 287.22223
$ws.Range("L22").Value = 720
$ws.Range("M22").Value = 62.77776999999998
$ws.Range("N22").Value = -1420

$ws.Range("H31").Value = 3691.9092
$ws.Range("I31").Value = 2104.2104
$ws.Range("J31").Value = 5846.643
$ws.Range("K31").Value = 2104.2104
$ws.Range("L31").Value = 5846.643
$ws.Range("M31").Value = -1809.2104
$ws.Range("N31").Value = -6436.643

$ws.Range("H34").Value = 3691.9092
$ws.Range("I34").Value = 2104.2104
$ws.Range("J34").Value = 5846.643
$ws.Range("K34").Value = 2104.2104
$ws.Range("L34").Value = 5846.643
$ws.Range("M34").Value = -1902.2104
$ws.Range("N34").Value = -6250.643

$ws.Range("H58").Value = 619292.4399999999
$ws.Range("I58").Value = 1176.5737
$ws.Range("J58").Value = 2504545.8
$ws.Range("K58").Value = 1176.5737
$ws.Range("L58").Value = 2504545.8
$ws.Range("M58").Value = -973.5736999999999
$ws.Range("N58").Value = -2504951.8

$ws.Range("H132").Value = 1992.48
$ws.Range("I132").Value = 1635.5
$ws.Range("J132").Value = 2322
$ws.Range("K132").Value = 4906.5
$ws.Range("L132").Value = 6966
$ws.Range("M132").Value = -2376.5
$ws.Range("N132").Value = -12026

$ws.Range("H134").Value = 2150.182
$ws.Range("I134").Value = 1171.7333
$ws.Range("J134").Value = 2965.5557
$ws.Range("K134").Value = 3515.199900000001
$ws.Range("L134").Value = 8896.667099999999
$ws.Range("M134").Value = -980.1999000000005
$ws.Range("N134").Value = -13966.6671

$ws.Range("H136").Value = 619292.4399999999
$ws.Range("I136").Value = 1176.5737
$ws.Range("J136").Value = 2504545.8
$ws.Range("K136").Value = 3529.7211
$ws.Range("L136").Value = 7513637.399999999
$ws.Range("M136").Value = -979.7210999999998
$ws.Range("N136").Value = -7518737.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 786.25
$ws.Range("I86").Value = 470
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1410
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = -224
$ws.Range("N86").Value = -11372

$ws.Range("H89").Value = 786.25
$ws.Range("I89").Value = 470
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 4230
$ws.Range("L89").Value = 27000
$ws.Range("M89").Value = 1698
$ws.Range("N89").Value = -38856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2049.6812
$ws.Range("I132").Value = 1375.2439
$ws.Range("J132").Value = 3037.25
$ws.Range("K132").Value = 4125.7317
$ws.Range("L132").Value = 9111.75
$ws.Range("M132").Value = -1595.7317
$ws.Range("N132").Value = -14171.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3739.8235
$ws.Range("I7").Value = 3449.1
$ws.Range("J7").Value = 4155.143
$ws.Range("K7").Value = 3449.1
$ws.Range("L7").Value = 4155.143
$ws.Range("M7").Value = -3337.1
$ws.Range("N7").Value = -4379.143

$ws.Range("H126").Value = 3739.8235
$ws.Range("I126").Value = 3449.1
$ws.Range("J126").Value = 4155.143
$ws.Range("K126").Value = 10347.3
$ws.Range("L126").Value = 12465.429
$ws.Range("M126").Value = -7877.299999999999
$ws.Range("N126").Value = -17405.429

$ws.Range("H132").Value = 18982.566
$ws.Range("I132").Value = 22361.535
$ws.Range("J132").Value = 9521.450000000001
$ws.Range("K132").Value = 67084.605
$ws.Range("L132").Value = 28564.35
$ws.Range("M132").Value = -64554.605
$ws.Range("N132").Value = -33624.35000000001

$ws.Range("H136").Value = 1674.1364
$ws.Range("I136").Value = 1225.9166
$ws.Range("J136").Value = 2212
$ws.Range("K136").Value = 3677.7498
$ws.Range("L136").Value = 6636
$ws.Range("M136").Value = -1127.7498
$ws.Range("N136").Value = -11736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 300
$ws.Range("I13").Value = 300
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -160

$ws.Range("H126").Value = 1955.4839
$ws.Range("I126").Value = 1708.2693
$ws.Range("J126").Value = 3241
$ws.Range("K126").Value = 5124.8079
$ws.Range("L126").Value = 9723
$ws.Range("M126").Value = -2654.8079
$ws.Range("N126").Value = -14663

$ws.Range("H132").Value = 1810.88
$ws.Range("I132").Value = 1732.3462
$ws.Range("J132").Value = 1895.9584
$ws.Range("K132").Value = 5197.0386
$ws.Range("L132").Value = 5687.8752
$ws.Range("M132").Value = -2667.0386
$ws.Range("N132").Value = -10747.8752

# Row 13's HQ profit column no longer applies (HQ price now 0), so the cell is cleared entirely
$ws.Range("N13").ClearContents()

Write-Output "done"
